$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 changes from "R40" to "1" (stored as text/shared string, not a number).
# Leading apostrophe forces Excel to keep the numeric-looking value as text
# instead of coercing it to a Number.
$ws.Range("B11").Value = "'1"
